$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Delete the C20 (10uF 1812) row entirely; everything below shifts up.
$ws.Rows("36").Delete()

# The former row 37 (1000pF/0603 capacitor) is now row 36.
# Its package cell (G36) switches from the numeric literal 603 to the text "0603".
$ws.Cells.Item(36, 7).NumberFormat = "@"
$ws.Cells.Item(36, 7).Value = "0603"
$ws.Cells.Item(36, 7).HorizontalAlignment = -4108

# Its value cell (E36) gets a dedicated font: Arial 10pt black, centered.
$ws.Cells.Item(36, 5).Font.Name = "Arial"
$ws.Cells.Item(36, 5).Font.Size = 10
$ws.Cells.Item(36, 5).Font.Color = 0
$ws.Cells.Item(36, 5).HorizontalAlignment = -4108

$ws.Range("E43").Select()
